$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "A2"
$ws.Range("B1").Value = "A30"
$ws.Range("A2").Value = "A17"
$ws.Range("B2").Value = "A20"
$ws.Range("A3").Value = "P1"
$ws.Range("B3").Value = "A20"

$ws.Range("D10").Select()
